# Adds quiz-2 / webassign-2 columns (G, H) to the roster sheet, mirroring
# the existing Q1/W1 (E/F) columns, and updates the sheet selection to the
# new H column to match the editor's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G1").Value = "Q2"
$ws.Range("H1").Value = "W2"

# Row 16 (Luna, Paloma) only earned 2 points on Q2; everyone else earned
# full credit (4) on both Q2 and W2.
$scores = @{16 = 2}

for ($r = 2; $r -le 33; $r++) {
    if ($scores.ContainsKey($r)) {
        $g = $scores[$r]
    } else {
        $g = 4
    }
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = 4
}

# Match the saved selection/active-cell state from the edit.
[void]$ws.Range("H2:H33").Select()
